$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The plate-layout grid (rows 3,5,7,9 = plate rows B,D,F,H) had its first
# data column (B) holding a stray/duplicate well value left over from a
# prior row. Fix each one so the interactive charts read the correct well
# id for that row - same value as used further down/up the column.
$ws.Range("B3").Value = "fae"
$ws.Range("B5").Value = "sea"
$ws.Range("B7").Value = "ase"
$ws.Range("B9").Value = "ca"

# Leave the active selection on B9, where the last edit was made.
[void]$ws.Range("B9").Select()
